# Auto-generated: apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.913.16"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.512.76"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'532.23"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").Value = "'138.50"
$ws.Range("E6").Value = "  -4.14%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("D9").Value = "2.514.93"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "'0.100"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "'0.354"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "2.956.76"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "'23.12"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("D16").Value = "58.871.34"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "2.515.22"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").Value = "'4.27"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "'322.54"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'5.80"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "'62.10"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "'0.424"
$ws.Range("E25").Value = "  -3.15%  "
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'6.70"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0768"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "'163.76"
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -9.58%  "
$ws.Range("D35").Value = "'1.44"
$ws.Range("E35").Value = "  -3.26%  "
$ws.Range("D36").Value = "'18.44"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("D37").Value = "'4.22"
$ws.Range("E37").Value = "  -3.22%  "
$ws.Range("D38").Value = "'1.57"
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'3.65"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "'0.803"
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").Value = "'5.20"
$ws.Range("E42").Value = "  -8.34%  "
$ws.Range("D43").Value = "'278.49"
$ws.Range("E43").Value = "  -6.53%  "
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "'0.0930"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").Value = "'18.39"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "'0.0510"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  -2.37%  "
